$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2, 2).Value = 33.7
$ws.Cells.Item(2, 3).Value = 13.236
$ws.Cells.Item(3, 2).Value = 83.5
$ws.Cells.Item(3, 3).Value = 4.983
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = 22.338
$ws.Cells.Item(5, 2).Value = 66
$ws.Cells.Item(5, 3).Value = 5.173
$ws.Cells.Item(6, 2).Value = 85.7
$ws.Cells.Item(6, 3).Value = 4.201
$ws.Cells.Item(7, 2).Value = 100
$ws.Cells.Item(7, 3).Value = 6.875
$ws.Cells.Item(8, 2).Value = 82.3
$ws.Cells.Item(8, 3).Value = 9.23
$ws.Cells.Item(9, 2).Value = 51.2
$ws.Cells.Item(9, 3).Value = 4.282
$ws.Cells.Item(10, 2).Value = 100
$ws.Cells.Item(10, 3).Value = 22
$ws.Cells.Item(11, 2).Value = 100
$ws.Cells.Item(11, 3).Value = 23
$ws.Cells.Item(12, 2).Value = 100
$ws.Cells.Item(12, 3).Value = 8.872999999999999
$ws.Cells.Item(13, 2).Value = 33.1
$ws.Cells.Item(13, 3).Value = 22.808
$ws.Cells.Item(14, 2).Value = 51.2
$ws.Cells.Item(14, 3).Value = 26.856
$ws.Cells.Item(15, 2).Value = 14.3
$ws.Cells.Item(15, 3).Value = 22.143
$ws.Cells.Item(16, 2).Value = 84.2
$ws.Cells.Item(16, 3).Value = 6.94
$ws.Cells.Item(17, 2).Value = 32.4
$ws.Cells.Item(17, 3).Value = 8.964
$ws.Cells.Item(18, 2).Value = 100
$ws.Cells.Item(18, 3).Value = 5.184
$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(19, 3).Value = 14.984
$ws.Cells.Item(20, 2).Value = 0
$ws.Cells.Item(20, 3).Value = 23
$ws.Cells.Item(21, 2).Value = 100
$ws.Cells.Item(21, 3).Value = 13.716
$ws.Cells.Item(22, 2).Value = 100
$ws.Cells.Item(22, 3).Value = 22.633
$ws.Cells.Item(23, 2).Value = 100
$ws.Cells.Item(23, 3).Value = 21
$ws.Cells.Item(24, 2).Value = 100
$ws.Cells.Item(24, 3).Value = 14
$ws.Cells.Item(25, 2).Value = 0
$ws.Cells.Item(25, 3).Value = 22
$ws.Cells.Item(26, 2).Value = 100
$ws.Cells.Item(26, 3).Value = 15.814
$ws.Cells.Item(27, 2).Value = 100
$ws.Cells.Item(27, 3).Value = 9.032
$ws.Cells.Item(28, 2).Value = 81.09999999999999
$ws.Cells.Item(28, 3).Value = 5.291
$ws.Cells.Item(29, 2).Value = 65
$ws.Cells.Item(29, 3).Value = 3.9
$ws.Cells.Item(30, 2).Value = 0
$ws.Cells.Item(30, 3).Value = 16.443
$ws.Cells.Item(31, 2).Value = 18.5
$ws.Cells.Item(31, 3).Value = 27.45
$ws.Cells.Item(32, 2).Value = 50.6
$ws.Cells.Item(32, 3).Value = 8.375
$ws.Cells.Item(33, 2).Value = 100
$ws.Cells.Item(33, 3).Value = 4.949
$ws.Cells.Item(34, 2).Value = 100
$ws.Cells.Item(34, 3).Value = 4.673
$ws.Cells.Item(35, 2).Value = 83.40000000000001
$ws.Cells.Item(35, 3).Value = 16.556
$ws.Cells.Item(36, 2).Value = 47.3
$ws.Cells.Item(36, 3).Value = 3.657
$ws.Cells.Item(37, 2).Value = 0
$ws.Cells.Item(37, 3).Value = 17.98
